# Swap the TC3 / TC4 test step content:
# TC3's "Steps" / "Expected Results" become what used to be TC4's,
# and TC4's "Steps" / "Expected Results" become what used to be TC3's.
# (The "TC3"/"TC4" id cells themselves stay where they are.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc3Steps    = $ws.Range("B25").Value2
$tc3Expected = $ws.Range("D25").Value2
$tc4Steps    = $ws.Range("B32").Value2
$tc4Expected = $ws.Range("D32").Value2

$ws.Range("B25").Value = $tc4Steps
$ws.Range("D25").Value = $tc4Expected
$ws.Range("B32").Value = $tc3Steps
$ws.Range("D32").Value = $tc3Expected
